$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 189 ("「忍耐は幸福の鍵」" post) - remaining rows shift up automatically,
# which also updates the sheet's used-range dimension from A1:C381 to A1:C380.
$ws.Rows.Item(189).Delete()
